# Auto-generated edit script applying numeric corrections to the
# profit-calculation columns (H-N) across several worksheets, per
# the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 4246.5386  # H62
$ws.Cells.Item(62, 9).Value = 4359  # I62
$ws.Cells.Item(62, 11).Value = 4359  # K62
$ws.Cells.Item(62, 13).Value = -3735  # M62

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 4246.5386  # H65
$ws.Cells.Item(65, 9).Value = 4359  # I65
$ws.Cells.Item(65, 11).Value = 21795  # K65
$ws.Cells.Item(65, 13).Value = -18675  # M65

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 0  # H127
$ws.Cells.Item(127, 9).Value = 0  # I127
$ws.Cells.Item(127, 11).Value = 0  # K127
$ws.Cells.Item(127, 13).ClearContents()  # M127

# ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value = 142935820  # H133
$ws.Cells.Item(133, 10).Value = 142935820  # J133
$ws.Cells.Item(133, 12).Value = 142935820  # L133
$ws.Cells.Item(133, 14).Value = -142945940  # N133

# ARM row 16
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(16, 8).Value = 3880  # H16
$ws.Cells.Item(16, 9).Value = 6324  # I16
$ws.Cells.Item(16, 10).Value = 825  # J16
$ws.Cells.Item(16, 11).Value = 6324  # K16
$ws.Cells.Item(16, 12).Value = 825  # L16
$ws.Cells.Item(16, 13).Value = -6037  # M16
$ws.Cells.Item(16, 14).Value = -1399  # N16

# ARM row 54
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(54, 8).Value = 82122  # H54
$ws.Cells.Item(54, 9).Value = 74245  # I54
$ws.Cells.Item(54, 11).Value = 74245  # K54
$ws.Cells.Item(54, 13).Value = -73476  # M54

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3018.4688  # H61
$ws.Cells.Item(61, 9).Value = 2992.8096  # I61
$ws.Cells.Item(61, 10).Value = 3067.4546  # J61
$ws.Cells.Item(61, 11).Value = 2992.8096  # K61
$ws.Cells.Item(61, 12).Value = 3067.4546  # L61
$ws.Cells.Item(61, 13).Value = -2780.8096  # M61
$ws.Cells.Item(61, 14).Value = -3491.4546  # N61

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 20856388  # H102
$ws.Cells.Item(102, 9).Value = 22729650  # I102
$ws.Cells.Item(102, 11).Value = 22729650  # K102
$ws.Cells.Item(102, 13).Value = -22728028  # M102

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 14302752  # H122
$ws.Cells.Item(122, 10).Value = 50733  # J122
$ws.Cells.Item(122, 12).Value = 152199  # L122
$ws.Cells.Item(122, 14).Value = -157099  # N122

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 3018.4688  # H136
$ws.Cells.Item(136, 9).Value = 2992.8096  # I136
$ws.Cells.Item(136, 10).Value = 3067.4546  # J136
$ws.Cells.Item(136, 11).Value = 8978.4288  # K136
$ws.Cells.Item(136, 12).Value = 9202.363799999999  # L136
$ws.Cells.Item(136, 13).Value = -6428.4288  # M136
$ws.Cells.Item(136, 14).Value = -14302.3638  # N136

# BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 0  # H35
$ws.Cells.Item(35, 9).Value = 0  # I35
$ws.Cells.Item(35, 11).Value = 0  # K35
$ws.Cells.Item(35, 13).ClearContents()  # M35

# BSM row 45
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(45, 8).Value = 0  # H45
$ws.Cells.Item(45, 9).Value = 0  # I45
$ws.Cells.Item(45, 11).Value = 0  # K45
$ws.Cells.Item(45, 13).ClearContents()  # M45

# BSM row 49
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(49, 8).Value = 19999  # H49
$ws.Cells.Item(49, 9).Value = 0  # I49
$ws.Cells.Item(49, 11).Value = 0  # K49
$ws.Cells.Item(49, 13).ClearContents()  # M49

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2209.4  # H86
$ws.Cells.Item(86, 9).Value = 2272.5833  # I86
$ws.Cells.Item(86, 10).Value = 1956.6666  # J86
$ws.Cells.Item(86, 11).Value = 2272.5833  # K86
$ws.Cells.Item(86, 12).Value = 1956.6666  # L86
$ws.Cells.Item(86, 13).Value = -1149.5833  # M86
$ws.Cells.Item(86, 14).Value = -4202.6666  # N86

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 2209.4  # H89
$ws.Cells.Item(89, 9).Value = 2272.5833  # I89
$ws.Cells.Item(89, 10).Value = 1956.6666  # J89
$ws.Cells.Item(89, 11).Value = 11362.9165  # K89
$ws.Cells.Item(89, 12).Value = 9783.333000000001  # L89
$ws.Cells.Item(89, 13).Value = -5746.916499999999  # M89
$ws.Cells.Item(89, 14).Value = -21015.333  # N89

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 28573264  # H94
$ws.Cells.Item(94, 9).Value = 1192.4286  # I94
$ws.Cells.Item(94, 11).Value = 1192.4286  # K94
$ws.Cells.Item(94, 13).Value = -741.4286  # M94

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 62502744  # H99
$ws.Cells.Item(99, 9).Value = 83335890  # I99
$ws.Cells.Item(99, 11).Value = 83335890  # K99
$ws.Cells.Item(99, 13).Value = -83334392  # M99

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 10418376  # H134
$ws.Cells.Item(134, 9).Value = 11365108  # I134
$ws.Cells.Item(134, 11).Value = 34095324  # K134
$ws.Cells.Item(134, 13).Value = -34092789  # M134

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 55559236  # H86
$ws.Cells.Item(86, 9).Value = 83336696  # I86
$ws.Cells.Item(86, 11).Value = 83336696  # K86
$ws.Cells.Item(86, 13).Value = -83335573  # M86

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 55559236  # H89
$ws.Cells.Item(89, 9).Value = 83336696  # I89
$ws.Cells.Item(89, 11).Value = 416683480  # K89
$ws.Cells.Item(89, 13).Value = -416677864  # M89

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 2200837.5  # H107
$ws.Cells.Item(107, 9).Value = 7857570.5  # I107
$ws.Cells.Item(107, 10).Value = 997.1389  # J107
$ws.Cells.Item(107, 11).Value = 7857570.5  # K107
$ws.Cells.Item(107, 12).Value = 997.1389  # L107
$ws.Cells.Item(107, 13).Value = -7855650.5  # M107
$ws.Cells.Item(107, 14).Value = -4837.1389  # N107

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 922.6667  # H134
$ws.Cells.Item(134, 9).Value = 881.2692  # I134
$ws.Cells.Item(134, 11).Value = 2643.8076  # K134
$ws.Cells.Item(134, 13).Value = -108.8076000000001  # M134

# CUL row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(126, 8).Value = 411  # H126
$ws.Cells.Item(126, 9).Value = 411  # I126
$ws.Cells.Item(126, 10).Value = 0  # J126
$ws.Cells.Item(126, 11).Value = 1233  # K126
$ws.Cells.Item(126, 13).Value = 3707  # M126
$ws.Cells.Item(126, 14).ClearContents()  # N126

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1696.138  # H131
$ws.Cells.Item(131, 9).Value = 1627.8182  # I131
$ws.Cells.Item(131, 10).Value = 1737.8889  # J131
$ws.Cells.Item(131, 11).Value = 4883.4546  # K131
$ws.Cells.Item(131, 12).Value = 5213.6667  # L131
$ws.Cells.Item(131, 13).Value = 156.5454  # M131
$ws.Cells.Item(131, 14).Value = -15293.6667  # N131

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 161.2  # H2
$ws.Cells.Item(2, 9).Value = 154.35294  # I2
$ws.Cells.Item(2, 11).Value = 154.35294  # K2
$ws.Cells.Item(2, 13).Value = -41.35293999999999  # M2

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 7522.4375  # H80
$ws.Cells.Item(80, 9).Value = 3626.5454  # I80
$ws.Cells.Item(80, 10).Value = 16093.4  # J80
$ws.Cells.Item(80, 11).Value = 3626.5454  # K80
$ws.Cells.Item(80, 12).Value = 16093.4  # L80
$ws.Cells.Item(80, 13).Value = -2628.5454  # M80
$ws.Cells.Item(80, 14).Value = -18089.4  # N80

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 7522.4375  # H83
$ws.Cells.Item(83, 9).Value = 3626.5454  # I83
$ws.Cells.Item(83, 10).Value = 16093.4  # J83
$ws.Cells.Item(83, 11).Value = 18132.727  # K83
$ws.Cells.Item(83, 12).Value = 80467  # L83
$ws.Cells.Item(83, 13).Value = -13140.727  # M83
$ws.Cells.Item(83, 14).Value = -90451  # N83

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 7239.6  # H22
$ws.Cells.Item(22, 9).Value = 1000  # I22
$ws.Cells.Item(22, 10).Value = 8799.5  # J22
$ws.Cells.Item(22, 11).Value = 1000  # K22
$ws.Cells.Item(22, 12).Value = 8799.5  # L22
$ws.Cells.Item(22, 13).Value = -705  # M22
$ws.Cells.Item(22, 14).Value = -9389.5  # N22

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 7239.6  # H27
$ws.Cells.Item(27, 9).Value = 1000  # I27
$ws.Cells.Item(27, 10).Value = 8799.5  # J27
$ws.Cells.Item(27, 11).Value = 1000  # K27
$ws.Cells.Item(27, 12).Value = 8799.5  # L27
$ws.Cells.Item(27, 13).Value = -893  # M27
$ws.Cells.Item(27, 14).Value = -9013.5  # N27

# LTW row 42
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 42569.6  # H42
$ws.Cells.Item(42, 9).Value = 57466  # I42
$ws.Cells.Item(42, 11).Value = 57466  # K42
$ws.Cells.Item(42, 13).Value = -56903  # M42

# LTW row 49
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(49, 8).Value = 42569.6  # H49
$ws.Cells.Item(49, 9).Value = 57466  # I49
$ws.Cells.Item(49, 11).Value = 57466  # K49
$ws.Cells.Item(49, 13).Value = -57319  # M49

# LTW row 50
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(50, 8).Value = 40056.75  # H50
$ws.Cells.Item(50, 9).Value = 40056.75  # I50
$ws.Cells.Item(50, 11).Value = 40056.75  # K50
$ws.Cells.Item(50, 13).Value = -39419.75  # M50

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1103.05  # H55
$ws.Cells.Item(55, 10).Value = 1481.5454  # J55
$ws.Cells.Item(55, 12).Value = 1481.5454  # L55
$ws.Cells.Item(55, 14).Value = -1827.5454  # N55

# WVR row 21
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value = 38311.25  # H21
$ws.Cells.Item(21, 9).Value = 38343  # I21
$ws.Cells.Item(21, 11).Value = 38343  # K21
$ws.Cells.Item(21, 13).Value = -38108  # M21

# WVR row 26
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(26, 8).Value = 60552.5  # H26
$ws.Cells.Item(26, 9).Value = 60552.5  # I26
$ws.Cells.Item(26, 11).Value = 60552.5  # K26
$ws.Cells.Item(26, 13).Value = -60259.5  # M26

# WVR row 35
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(35, 8).Value = 38311.25  # H35
$ws.Cells.Item(35, 9).Value = 38343  # I35
$ws.Cells.Item(35, 11).Value = 38343  # K35
$ws.Cells.Item(35, 13).Value = -38053  # M35

# WVR row 37
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(37, 8).Value = 37499  # H37
$ws.Cells.Item(37, 9).Value = 37499  # I37
$ws.Cells.Item(37, 11).Value = 37499  # K37
$ws.Cells.Item(37, 13).Value = -37296  # M37

# WVR row 49
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 39999  # H49
$ws.Cells.Item(49, 9).Value = 39999  # I49
$ws.Cells.Item(49, 11).Value = 39999  # K49
$ws.Cells.Item(49, 13).Value = -39769  # M49

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 32735.438  # H54
$ws.Cells.Item(54, 9).Value = 29384.5  # I54
$ws.Cells.Item(54, 10).Value = 33214.145  # J54
$ws.Cells.Item(54, 11).Value = 29384.5  # K54
$ws.Cells.Item(54, 12).Value = 33214.145  # L54
$ws.Cells.Item(54, 13).Value = -28864.5  # M54
$ws.Cells.Item(54, 14).Value = -34254.145  # N54

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 9260826  # H96
$ws.Cells.Item(96, 10).Value = 1884  # J96
$ws.Cells.Item(96, 12).Value = 1884  # L96
$ws.Cells.Item(96, 14).Value = -4630  # N96
